$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - RIOT
$ws.Range("K2").Value = 59.6
$ws.Range("N2").Value = 54.77309453746771

# Row 3 - BTC-USD
$ws.Range("D3").Value = 92123.83
$ws.Range("E3").Value = 65.59999999999999
$ws.Range("F3").Value = 1.91
$ws.Range("H3").Value = 53
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 55.6
$ws.Range("N3").Value = 54.77309453746771

# Row 4 - COIN
$ws.Range("K4").Value = 51.4
$ws.Range("N4").Value = 54.77309453746771

# Row 5 - MARA
$ws.Range("K5").Value = 49.6
$ws.Range("N5").Value = 54.77309453746771

# Row 6 - MSTR
$ws.Range("K6").Value = 35.8
$ws.Range("N6").Value = 54.77309453746771
